$d = $word.ActiveDocument

# 1) Append " - V1" (bold) right after the title text, as a separate run.
$r = $d.Content
$r.Find.Execute("Sistema de Processamento e Análise de Imagens", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)  # wdCollapseEnd
$r.Font.Bold = $true
$r.InsertAfter(" – V1")

# 2) Mark the "Fontepargpadro" (Default Paragraph Font) character style as SemiHidden.
$style = $d.Styles("Fontepargpadro")
$style.Visibility = $false
